$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "Value rounded"
$ws.Range("D1").Value = "Value scientific notation"
$ws.Range("F1").Value = "Distance (m)"
$ws.Range("H1").Value = "Elevation (m)"
$ws.Range("I1").Value = "Hill height (m)"
$ws.Range("J1").Value = "FIPs"
$ws.Range("L1").Value = "UTM easting"
$ws.Range("M1").Value = "UTM northing"
$ws.Range("P1").Value = "Receptor type"
